$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 933.6923
$ws.Range("I38").Value = 143.11111
$ws.Range("J38").Value = 2712.5
$ws.Range("K38").Value = 429.33333
$ws.Range("L38").Value = 8137.5
$ws.Range("M38").Value = -57.33332999999999
$ws.Range("N38").Value = -8881.5

# Row 40
$ws.Range("H40").Value = 2163.7
$ws.Range("J40").Value = 2163.7
$ws.Range("L40").Value = 2163.7
$ws.Range("N40").Value = -2513.7

# Row 58
$ws.Range("H58").Value = 3468.4211
$ws.Range("I58").Value = 533.3333
$ws.Range("J58").Value = 4823.077
$ws.Range("K58").Value = 1599.9999
$ws.Range("L58").Value = 14469.231
$ws.Range("M58").Value = -1449.9999
$ws.Range("N58").Value = -14769.231

# Row 86
$ws.Range("H86").Value = 3053.6296
$ws.Range("I86").Value = 2066
$ws.Range("J86").Value = 4490.1816
$ws.Range("K86").Value = 2066
$ws.Range("L86").Value = 4490.1816
$ws.Range("M86").Value = -943
$ws.Range("N86").Value = -6736.1816

# Row 89
$ws.Range("H89").Value = 3053.6296
$ws.Range("I89").Value = 2066
$ws.Range("J89").Value = 4490.1816
$ws.Range("K89").Value = 10330
$ws.Range("L89").Value = 22450.908
$ws.Range("M89").Value = -4714
$ws.Range("N89").Value = -33682.908

# Row 107
$ws.Range("H107").Value = 1577.8572
$ws.Range("I107").Value = 1325.2667
$ws.Range("J107").Value = 2209.3333
$ws.Range("K107").Value = 1325.2667
$ws.Range("L107").Value = 2209.3333
$ws.Range("M107").Value = 594.7333000000001
$ws.Range("N107").Value = -6049.3333

# Row 132
$ws.Range("H132").Value = 10424081
$ws.Range("I132").Value = 12823253
$ws.Range("K132").Value = 38469759
$ws.Range("M132").Value = -38467229

# Row 141
$ws.Range("H141").Value = 1100
$ws.Range("I141").Value = 1095
$ws.Range("J141").Value = 1105
$ws.Range("K141").Value = 3285
$ws.Range("L141").Value = 3315
$ws.Range("M141").Value = 1895
$ws.Range("N141").Value = -13675


# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 12063.444
$ws.Range("I2").Value = 1109
$ws.Range("J2").Value = 25756.5
$ws.Range("K2").Value = 1109
$ws.Range("L2").Value = 25756.5
$ws.Range("M2").Value = -996
$ws.Range("N2").Value = -25982.5

# Row 5
$ws.Range("H5").Value = 223.5
$ws.Range("I5").Value = 135.25
$ws.Range("K5").Value = 135.25
$ws.Range("M5").Value = -23.25

# Row 32
$ws.Range("H32").Value = 3994.795
$ws.Range("I32").Value = 4474.788
$ws.Range("K32").Value = 4474.788
$ws.Range("M32").Value = -4187.788

# Row 45
$ws.Range("H45").Value = 1215.1666
$ws.Range("I45").Value = 1450.8334
$ws.Range("J45").Value = 743.8333
$ws.Range("K45").Value = 1450.8334
$ws.Range("L45").Value = 743.8333
$ws.Range("M45").Value = -1073.8334
$ws.Range("N45").Value = -1497.8333

# Row 97
$ws.Range("H97").Value = 714.0833
$ws.Range("I97").Value = 371.25
$ws.Range("J97").Value = 1399.75
$ws.Range("K97").Value = 371.25
$ws.Range("L97").Value = 1399.75
$ws.Range("M97").Value = 124.75
$ws.Range("N97").Value = -2391.75

# Row 102
$ws.Range("H102").Value = 9804795
$ws.Range("I102").Value = 10417470
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 10417470
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -10415848
$ws.Range("N102").Value = -5244

# Row 116
$ws.Range("H116").Value = 12063.444
$ws.Range("I116").Value = 1109
$ws.Range("J116").Value = 25756.5
$ws.Range("K116").Value = 1109
$ws.Range("L116").Value = 25756.5
$ws.Range("M116").Value = 1185
$ws.Range("N116").Value = -30344.5

# Row 122
$ws.Range("H122").Value = 1505.091
$ws.Range("I122").Value = 1358.7646
$ws.Range("J122").Value = 2002.6
$ws.Range("K122").Value = 4076.2938
$ws.Range("L122").Value = 6007.799999999999
$ws.Range("M122").Value = -1626.2938
$ws.Range("N122").Value = -10907.8

# Row 132
$ws.Range("H132").Value = 2771.125
$ws.Range("I132").Value = 2351.8823
$ws.Range("J132").Value = 3789.2856
$ws.Range("K132").Value = 7055.646900000001
$ws.Range("L132").Value = 11367.8568
$ws.Range("M132").Value = -4525.646900000001
$ws.Range("N132").Value = -16427.8568


# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 12063.444
$ws.Range("I3").Value = 1109
$ws.Range("J3").Value = 25756.5
$ws.Range("K3").Value = 1109
$ws.Range("L3").Value = 25756.5
$ws.Range("M3").Value = -995
$ws.Range("N3").Value = -25984.5

# Row 4
$ws.Range("H4").Value = 223.5
$ws.Range("I4").Value = 135.25
$ws.Range("K4").Value = 135.25
$ws.Range("M4").Value = -20.25

# Row 94
$ws.Range("H94").Value = 14706466
$ws.Range("I94").Value = 15625495
$ws.Range("K94").Value = 15625495
$ws.Range("M94").Value = -15625044


# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 125001250
$ws.Range("I16").Value = 142858350
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 142858350
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -142858063
$ws.Range("N16").Value = -2074

# Row 22
$ws.Range("H22").Value = 541.25
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 582.5
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 582.5
$ws.Range("M22").Value = -150
$ws.Range("N22").Value = -1282.5

# Row 31
$ws.Range("H31").Value = 1854.561
$ws.Range("I31").Value = 833.6842
$ws.Range("J31").Value = 2736.2273
$ws.Range("K31").Value = 833.6842
$ws.Range("L31").Value = 2736.2273
$ws.Range("M31").Value = -538.6842
$ws.Range("N31").Value = -3326.2273

# Row 34
$ws.Range("H34").Value = 1854.561
$ws.Range("I34").Value = 833.6842
$ws.Range("J34").Value = 2736.2273
$ws.Range("K34").Value = 833.6842
$ws.Range("L34").Value = 2736.2273
$ws.Range("M34").Value = -631.6842
$ws.Range("N34").Value = -3140.2273

# Row 62
$ws.Range("H62").Value = 7145778
$ws.Range("I62").Value = 3029.074
$ws.Range("K62").Value = 3029.074
$ws.Range("M62").Value = -2405.074

# Row 65
$ws.Range("H65").Value = 7145778
$ws.Range("I65").Value = 3029.074
$ws.Range("K65").Value = 15145.37
$ws.Range("M65").Value = -12025.37

# Row 113
$ws.Range("H113").Value = 125001250
$ws.Range("I113").Value = 142858350
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 142858350
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = -142856180
$ws.Range("N113").Value = -5840

# Row 132
$ws.Range("H132").Value = 2432.5715
$ws.Range("I132").Value = 1294.75
$ws.Range("J132").Value = 3949.6667
$ws.Range("K132").Value = 3884.25
$ws.Range("L132").Value = 11849.0001
$ws.Range("M132").Value = -1354.25
$ws.Range("N132").Value = -16909.0001

# Row 134
$ws.Range("H134").Value = 27779380
$ws.Range("I134").Value = 33334756
$ws.Range("K134").Value = 100004268
$ws.Range("M134").Value = -100001733


# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1337.7916
$ws.Range("I5").Value = 1407.1364
$ws.Range("K5").Value = 4221.4092
$ws.Range("M5").Value = -4109.4092

# Row 131
$ws.Range("H131").Value = 14926636
$ws.Range("J131").Value = 1353.6721
$ws.Range("L131").Value = 4061.0163
$ws.Range("N131").Value = -14141.0163

# Row 135
$ws.Range("H135").Value = 1337.7916
$ws.Range("I135").Value = 1407.1364
$ws.Range("K135").Value = 12664.2276
$ws.Range("M135").Value = -10129.2276

# Row 139
$ws.Range("H139").Value = 2578.7646
$ws.Range("I139").Value = 2207.25
$ws.Range("K139").Value = 6621.75
$ws.Range("M139").Value = -1481.75

# Row 140
$ws.Range("H140").Value = 37254.324
$ws.Range("I140").Value = 47266.582
$ws.Range("K140").Value = 141799.746
$ws.Range("M140").Value = -136619.746


# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 4367.6665
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 4367.6665
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 4367.6665
$ws.Range("N80").Value = -6363.6665
$ws.Range("M80").ClearContents()

# Row 83
$ws.Range("H83").Value = 4367.6665
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 4367.6665
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 21838.3325
$ws.Range("N83").Value = -31822.3325
$ws.Range("M83").ClearContents()


# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1564.9166
$ws.Range("I22").Value = 1422.2222
$ws.Range("K22").Value = 1422.2222
$ws.Range("M22").Value = -1127.2222

# Row 27
$ws.Range("H27").Value = 1564.9166
$ws.Range("I27").Value = 1422.2222
$ws.Range("K27").Value = 1422.2222
$ws.Range("M27").Value = -1315.2222


# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 3791.1177
$ws.Range("I132").Value = 3759.8333
$ws.Range("J132").Value = 3866.2
$ws.Range("K132").Value = 11279.4999
$ws.Range("L132").Value = 11598.6
$ws.Range("M132").Value = -8749.499899999999
$ws.Range("N132").Value = -16658.6


